# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (D2) and "Correspond Handback DateTime" (G2)
# for the first file (10fb5430-...) in both the zh-cn and de-de handback-status sheets,
# reflecting a newly generated handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-25 13:08:46"
$wsZhCn.Range("G2").Value = "2016-01-25 13:09:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-25 13:08:56"
$wsDeDe.Range("G2").Value = "2016-01-25 13:09:46"
